# Generate Report for Handback
#
# The row for file "31137417-6e4e-4f30-a362-a0e4ca2830f2" failed its
# handback transform (the file name returned by the translation vendor
# did not match the file name that was handed off). Reflect that in the
# per-language status sheets:
#   - Status column (C3) moves from "Ready for handoff" to
#     "Handback transform failed"
#   - A new "Error Detail" (column L) explains the mismatch for each
#     language.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Update the status for the failed handback. The same text is used on
# the Overview sheet (columns B and C) as well as on each per-language
# sheet's Status column (C) since they all describe the same file.
$wsOverview.Range("B3").Value = "Handback transform failed"
$wsOverview.Range("C3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# Record the error detail explaining the handback failure for each
# target language.
$wsZhCn.Range("L3").Value = "Handback file name: uhmstdki.l25 is different with handoff file name: 31137417-6e4e-4f30-a362-a0e4ca2830f2.388609cd928aa99ec2c8b66238cd8107b7499b4d.zh-cn."
$wsDeDe.Range("L3").Value = "Handback file name: uhmstdki.l25 is different with handoff file name: 31137417-6e4e-4f30-a362-a0e4ca2830f2.388609cd928aa99ec2c8b66238cd8107b7499b4d.de-de."
